$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must remain text (matches original
# inline-string cell type / display, e.g. '1.00', '0.587') - force text format
# first so Excel doesn't coerce the literal into a Double and mangle the
# formatting/precision (e.g. '1.00' -> 1, '0.587' -> 0.58699999999999997).
$textCells = @("D5", "D6", "D7", "D8", "D10", "D11", "D15", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D44", "D45", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell-by-cell, in sheet order.
$ws.Range("D2").Value = '59.254.56'
$ws.Range("E2").Value = '  +4.89%  '
$ws.Range("D3").Value = '3.346.81'
$ws.Range("E3").Value = '  +2.81%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '410.30'
$ws.Range("E5").Value = '  +3.09%  '
$ws.Range("D6").Value = '111.65'
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("D7").Value = '0.587'
$ws.Range("E7").Value = '  +4.73%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  +2.09%  '
$ws.Range("D10").Value = '40.19'
$ws.Range("E10").Value = '  +1.99%  '
$ws.Range("D11").Value = '0.0986'
$ws.Range("E11").Value = '  +2.81%  '
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("D13").Value = '3.874.19'
$ws.Range("E13").Value = '  +3.04%  '
$ws.Range("E14").Value = '  +4.34%  '
$ws.Range("D15").Value = '19.42'
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").Value = '3.333.88'
$ws.Range("E16").Value = '  +2.35%  '
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").Value = '59.036.44'
$ws.Range("E18").Value = '  +4.58%  '
$ws.Range("D19").Value = '10.84'
$ws.Range("E19").Value = '  -1.69%  '
$ws.Range("D20").Value = '3.36'
$ws.Range("E20").Value = '  +1.00%  '
$ws.Range("D21").Value = '0.0000112'
$ws.Range("E21").Value = '  +6.97%  '
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").Value = '304.70'
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("D24").Value = '75.63'
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("D25").Value = '3.20'
$ws.Range("E25").Value = '  -1.17%  '
$ws.Range("D26").Value = '28.67'
$ws.Range("E26").Value = '  +1.43%  '
$ws.Range("D27").Value = '4.47'
$ws.Range("E27").Value = '  +2.79%  '
$ws.Range("D28").Value = '7.89'
$ws.Range("E28").Value = '  -3.17%  '
$ws.Range("E29").Value = '  +1.92%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '0.118'
$ws.Range("E30").Value = '  +5.89%  '
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").Value = '0.172'
$ws.Range("E31").Value = '  +1.16%  '
$ws.Range("D32").Value = '11.65'
$ws.Range("E32").Value = '  +5.00%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").Value = '40.39'
$ws.Range("E34").Value = '  +9.45%  '
$ws.Range("D35").Value = '0.0524'
$ws.Range("E35").Value = '  +7.45%  '
$ws.Range("D36").Value = '2.13'
$ws.Range("E36").Value = '  +0.62%  '
$ws.Range("D37").Value = '51.99'
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("D38").Value = '3.11'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").Value = '3.44'
$ws.Range("E40").Value = '  -2.62%  '
$ws.Range("D41").Value = '137.39'
$ws.Range("E41").Value = '  +2.29%  '
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("D44").Value = '3.98'
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("D45").Value = '16.96'
$ws.Range("E45").Value = '  -3.89%  '
$ws.Range("D46").Value = '0.278'
$ws.Range("E46").Value = '  -2.56%  '
$ws.Range("D47").Value = '2.27'
$ws.Range("E47").Value = '  +8.65%  '
$ws.Range("D48").Value = '22.45'
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("D49").Value = '2.205.54'
$ws.Range("E49").Value = '  +2.56%  '
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").Value = '6.39'
$ws.Range("E51").Value = '  +6.47%  '
